$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 responses: Q2_2 and Q2_3 change to "Excellent 5",
# Q3_3 and Q3_2 change to "Good 4"
$ws.Range("G3").Value = "Excellent 5"
$ws.Range("H3").Value = "Excellent 5"
$ws.Range("J3").Value = "Good 4"
$ws.Range("K3").Value = "Good 4"

# Update the active selection to K3
$ws.Range("K3").Select()
